# Apply updated crypto market data (price/volume) scraped on 2023-01-05.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D (Price) and E (Volume(1h)) are stored as text in the sheet,
# so a leading apostrophe forces Excel to keep them as text instead of numbers.
$ws.Range("D2").Value = "'256.46"
$ws.Range("E2").Value = "'-0.48%"
$ws.Range("D3").Value = "'27.06"
$ws.Range("E3").Value = "'-3.50%"
$ws.Range("D4").Value = "'4.656"
$ws.Range("E4").Value = "'-10.86%"
$ws.Range("D5").Value = "'0.05892"
$ws.Range("E5").Value = "'-0.23%"
$ws.Range("D6").Value = "'6.634"
$ws.Range("E6").Value = "'-1.02%"
$ws.Range("D7").Value = "'0.8626"
$ws.Range("E7").Value = "'-0.74%"
$ws.Range("D8").Value = "'0.9309"
$ws.Range("E8").Value = "'-10.49%"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01039"
$ws.Range("E9").Value = "'1,614.81%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1408"
$ws.Range("E10").Value = "'-0.31%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03715"
$ws.Range("E11").Value = "'2.23%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07090"
$ws.Range("E12").Value = "'-1.35%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03231"
$ws.Range("E13").Value = "'2.62%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09218"
$ws.Range("E14").Value = "'0.01%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001547"
$ws.Range("E15").Value = "'0.62%"
$ws.Range("D16").Value = "'0.006077"
$ws.Range("E16").Value = "'3.31%"
$ws.Range("E17").Value = "'0.41%"
$ws.Range("D18").Value = "'3.192"
$ws.Range("E18").Value = "'-1.11%"
$ws.Range("D19").Value = "'2.202"
$ws.Range("E19").Value = "'-1.08%"
$ws.Range("D20").Value = "'0.3100"
$ws.Range("E20").Value = "'-0.68%"
$ws.Range("D21").Value = "'0.1275"
$ws.Range("E21").Value = "'-1.12%"
$ws.Range("D22").Value = "'3.848"
$ws.Range("E22").Value = "'9.13%"
$ws.Range("D23").Value = "'0.04232"
$ws.Range("E23").Value = "'0.80%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'0.17%"
$ws.Range("D25").Value = "'0.004278"
$ws.Range("E27").Value = "'31.63%"
$ws.Range("D40").Value = "'0.03824"
$ws.Range("E40").Value = "'-0.29%"
$ws.Range("D41").Value = "'0.006221"
$ws.Range("E41").Value = "'14.84%"
$ws.Range("E42").Value = "'-0.47%"
$ws.Range("E43").Value = "'-4.29%"
$ws.Range("D44").Value = "'0.01137"
$ws.Range("E44").Value = "'6.63%"
$ws.Range("D45").Value = "'0.00005460"
$ws.Range("E45").Value = "'0.86%"
$ws.Range("E46").Value = "'0.06%"
$ws.Range("D48").Value = "'0.002279"
$ws.Range("E48").Value = "'6.69%"
$ws.Range("E49").Value = "'0.06%"
$ws.Range("E50").Value = "'0.06%"
